$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The seven cells on Sheet1 hold SQL query text (as shared strings) that all
# repeat the same LEFT JOIN block. The join columns changed from the
# generic "id" to the fully-qualified "study_id" / "participant_id" names.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellRef in $cells) {
    $range = $ws.Range($cellRef)
    $text = $range.Value2

    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

    $range.Value2 = $text
}

# Column C was widened from an auto (best-fit) width to a fixed, wider width.
$ws.Columns("C").ColumnWidth = 68.25

# The active cell / scroll position moved down to the last query (row 7).
$ws.Range("C7").Select()
